$wb = $excel.ActiveWorkbook

# Sheet "OFF" - row 3 (R) updates
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 403
$wsOff.Range("C3").Value = 287
$wsOff.Range("D3").Value = 81
$wsOff.Range("E3").Value = 39
$wsOff.Range("G3").Value = 8

# Sheet "DEF" - row 3 (R) updates
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 334
$wsDef.Range("C3").Value = 229
$wsDef.Range("D3").Value = 96
$wsDef.Range("E3").Value = 40
$wsDef.Range("G3").Value = 7
